$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "25.568.95"
$ws.Cells.Item(2,5).Value = "  +2.44%  "

$ws.Cells.Item(3,4).Value = "1.671.14"
$ws.Cells.Item(3,5).Value = "  +1.85%  "

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "0.9988"
$ws.Cells.Item(4,4).ClearFormats()
$ws.Cells.Item(4,5).Value = "  +0.01%  "

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "238.99"
$ws.Cells.Item(5,4).ClearFormats()
$ws.Cells.Item(5,5).Value = "  +1.48%  "

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "0.9997"
$ws.Cells.Item(6,4).ClearFormats()
$ws.Cells.Item(6,5).Value = "  -0.13%  "

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "0.4811"
$ws.Cells.Item(7,4).ClearFormats()
$ws.Cells.Item(7,5).Value = "  +0.68%  "

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.2632"
$ws.Cells.Item(8,4).ClearFormats()
$ws.Cells.Item(8,5).Value = "  +2.38%  "

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "0.06187"
$ws.Cells.Item(9,4).ClearFormats()
$ws.Cells.Item(9,5).Value = "  +3.22%  "

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.07016"
$ws.Cells.Item(10,4).ClearFormats()
$ws.Cells.Item(10,5).Value = "  -2.72%  "

$ws.Cells.Item(11,4).Value = "1.670.06"
$ws.Cells.Item(11,5).Value = "  +1.78%  "

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "14.92"
$ws.Cells.Item(12,4).ClearFormats()
$ws.Cells.Item(12,5).Value = "  +0.87%  "

$ws.Cells.Item(13,5).Value = "  -3.88%  "

$ws.Cells.Item(14,5).Value = "  -2.18%  "

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "75.35"
$ws.Cells.Item(15,4).ClearFormats()
$ws.Cells.Item(15,5).Value = "  +3.71%  "

$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "0.9999"
$ws.Cells.Item(16,4).ClearFormats()
$ws.Cells.Item(16,5).Value = "  -0.15%  "

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "0.9994"
$ws.Cells.Item(17,4).ClearFormats()
$ws.Cells.Item(17,5).Value = "  +0.04%  "

$ws.Cells.Item(18,4).Value = "25.566.29"
$ws.Cells.Item(18,5).Value = "  +2.50%  "

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "0.000006786"
$ws.Cells.Item(19,4).ClearFormats()
$ws.Cells.Item(19,5).Value = "  +2.92%  "

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "11.47"
$ws.Cells.Item(20,4).ClearFormats()
$ws.Cells.Item(20,5).Value = "  +1.57%  "

$ws.Cells.Item(21,4).Value = "1.882.77"
$ws.Cells.Item(21,5).Value = "  +1.59%  "

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "4.453"
$ws.Cells.Item(22,4).ClearFormats()
$ws.Cells.Item(22,5).Value = "  -0.10%  "

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "8.740"
$ws.Cells.Item(23,4).ClearFormats()
$ws.Cells.Item(23,5).Value = "  +1.89%  "

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "5.293"
$ws.Cells.Item(24,4).ClearFormats()
$ws.Cells.Item(24,5).Value = "  +0.37%  "

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "136.89"
$ws.Cells.Item(25,4).ClearFormats()
$ws.Cells.Item(25,5).Value = "  +3.63%  "

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "15.08"
$ws.Cells.Item(26,4).ClearFormats()
$ws.Cells.Item(26,5).Value = "  +1.76%  "

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "1.393"
$ws.Cells.Item(27,4).ClearFormats()
$ws.Cells.Item(27,5).Value = "  +0.80%  "

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "1.731"
$ws.Cells.Item(28,4).ClearFormats()
$ws.Cells.Item(28,5).Value = "  +4.35%  "

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "104.93"
$ws.Cells.Item(29,4).ClearFormats()
$ws.Cells.Item(29,5).Value = "  +1.82%  "

$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "3.996"
$ws.Cells.Item(30,4).ClearFormats()
$ws.Cells.Item(30,5).Value = "  +7.31%  "

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "0.07817"
$ws.Cells.Item(31,4).ClearFormats()
$ws.Cells.Item(31,5).Value = "  +0.24%  "

$ws.Cells.Item(32,5).Value = "  +3.45%  "

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "0.9987"
$ws.Cells.Item(33,4).ClearFormats()
$ws.Cells.Item(33,5).Value = "  -0.11%  "

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "0.04238"
$ws.Cells.Item(34,4).ClearFormats()
$ws.Cells.Item(34,5).Value = "  -3.87%  "

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "2.616"
$ws.Cells.Item(35,4).ClearFormats()
$ws.Cells.Item(35,5).Value = "  +0.96%  "

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "0.6109"
$ws.Cells.Item(36,4).ClearFormats()
$ws.Cells.Item(36,5).Value = "  +4.76%  "

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = "0.9517"
$ws.Cells.Item(37,4).ClearFormats()
$ws.Cells.Item(37,5).Value = "  +2.87%  "

$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "2.597"
$ws.Cells.Item(38,4).ClearFormats()
$ws.Cells.Item(38,5).Value = "  +1.68%  "

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "0.8596"
$ws.Cells.Item(39,4).ClearFormats()
$ws.Cells.Item(39,5).Value = "  +2.15%  "

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "0.9993"
$ws.Cells.Item(40,4).ClearFormats()
$ws.Cells.Item(40,5).Value = "  +0.01%  "

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "1.866"
$ws.Cells.Item(41,4).ClearFormats()
$ws.Cells.Item(41,5).Value = "  +3.80%  "

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "0.01475"
$ws.Cells.Item(42,4).ClearFormats()
$ws.Cells.Item(42,5).Value = "  -5.53%  "

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "96.21"
$ws.Cells.Item(43,4).ClearFormats()
$ws.Cells.Item(43,5).Value = "  -1.03%  "

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "0.3779"
$ws.Cells.Item(44,4).ClearFormats()
$ws.Cells.Item(44,5).Value = "  +1.85%  "

$ws.Cells.Item(46,5).Value = "  -2.58%  "

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "6.227"
$ws.Cells.Item(47,4).ClearFormats()
$ws.Cells.Item(47,5).Value = "  +2.43%  "

$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "0.05256"
$ws.Cells.Item(48,4).ClearFormats()

$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "29.88"
$ws.Cells.Item(49,4).ClearFormats()
$ws.Cells.Item(49,5).Value = "  +0.85%  "

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "7.365"
$ws.Cells.Item(50,4).ClearFormats()
$ws.Cells.Item(50,5).Value = "  +2.25%  "

$ws.Cells.Item(51,5).Value = "  +0.19%  "
